$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.925.09"
$ws.Range("E2").Value = "  -3.56%  "

$ws.Range("D3").Value = "'3.530.75"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'609.81"
$ws.Range("E5").Value = "  -5.34%  "

$ws.Range("D6").Value = "'152.52"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("D7").Value = "'3.528.78"
$ws.Range("E7").Value = "  -3.75%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "'0.483"
$ws.Range("E9").Value = "  -3.06%  "

$ws.Range("E10").Value = "  -3.35%  "

$ws.Range("D11").Value = "'6.88"
$ws.Range("E11").Value = "  -2.74%  "

$ws.Range("D12").Value = "'0.426"
$ws.Range("E12").Value = "  -3.28%  "

$ws.Range("D13").Value = "'0.0000220"
$ws.Range("E13").Value = "  -4.36%  "

$ws.Range("D14").Value = "'4.128.94"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").Value = "'31.67"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("D16").Value = "'3.524.96"
$ws.Range("E16").Value = "  -4.11%  "

$ws.Range("D17").Value = "'66.949.30"
$ws.Range("E17").Value = "  -3.48%  "

$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.28"
$ws.Range("E19").Value = "  -2.63%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'15.36"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("D21").Value = "'443.92"
$ws.Range("E21").Value = "  -4.88%  "

$ws.Range("D22").Value = "'9.20"
$ws.Range("E22").Value = "  -8.16%  "

$ws.Range("D23").Value = "'0.630"
$ws.Range("E23").Value = "  -2.34%  "

$ws.Range("D24").Value = "'77.52"
$ws.Range("E24").Value = "  -2.15%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "'3.669.13"
$ws.Range("E25").Value = "  -3.74%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "'0.0000123"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = "  -6.20%  "

$ws.Range("D29").Value = "'8.19"
$ws.Range("E29").Value = "  -9.02%  "

$ws.Range("D30").Value = "'2.53"
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").Value = "'1.67"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "'25.73"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("D34").Value = "'0.158"
$ws.Range("E34").Value = "  -2.97%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.14"
$ws.Range("E35").Value = "  -4.15%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.87"
$ws.Range("E36").Value = "  -6.35%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "'3.524.39"
$ws.Range("E37").Value = "  -3.82%  "

$ws.Range("D38").Value = "'8.02"
$ws.Range("E38").Value = "  -4.73%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").Value = "'174.06"
$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'5.55"
$ws.Range("E42").Value = "  -5.14%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.12"
$ws.Range("E43").Value = "  -3.13%  "

$ws.Range("D44").Value = "'0.0859"
$ws.Range("E44").Value = "  -3.67%  "

$ws.Range("D45").Value = "'0.890"
$ws.Range("E45").Value = "  -3.48%  "

$ws.Range("D46").Value = "'45.31"
$ws.Range("E46").Value = "  -3.78%  "

$ws.Range("D47").Value = "'27.16"
$ws.Range("E47").Value = "  -5.43%  "

$ws.Range("D48").Value = "'2.55"
$ws.Range("E48").Value = "  -4.50%  "

$ws.Range("D49").Value = "'1.22"
$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").Value = "'7.55"
$ws.Range("E50").Value = "  -2.93%  "

$ws.Range("E51").Value = "  -3.26%  "
